# Apply the edit described by the diff:
#  - Remove column F entirely (the "pval" header and per-row p-value cells)
#  - Remove the "Int. vs. Unint. p-val: 1.796443e-12 ***" row (old row 8
#    text), replacing row 8's note/SD-in-parens text from the row below it
#    (old row 9), then delete that now-redundant row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 keeps its "7" label (column A) but its B/C text is replaced with
# what used to live in row 9 (the scale note + "SD in parens." caption).
$ws.Range("B8").Value = "Note: Scale was 1 (not linked) to 7 (tightly linked)."
$ws.Range("C8").Value = "SD in parens."

# Row 9 is now redundant - remove it entirely.
$ws.Range("A9:E9").EntireRow.Delete()

# Delete the whole column F (pval column).
$ws.Range("F1:F8").EntireColumn.Delete()
